$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$np = $s.NotesPage
try {
  $ph = $np.Shapes.AddPlaceholder(2)
  Write-Host "Placeholder added: $($ph.Name)"
  $ph.TextFrame.TextRange.Text = "Now I'm going to explain you how RL works really. Giving the past example with the cat, imagine you have the cat here as the agent"
} catch {
  Write-Host "Error: $_"
}
